$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Insert a new "Read me" worksheet at the very front of the
#    workbook and drop an explanatory textbox (shape) on it.
# ------------------------------------------------------------------
$firstSheet = $wb.Worksheets.Item(1)
$readme = $wb.Worksheets.Add($firstSheet)
$readme.Name = "Read me"

$shp = $readme.Shapes.AddTextbox(1, 0, 0, 639, 224)
$shp.Name = "TextBox 1"

$readmeText = "Binary predictions of unobserved ('unknown') hosts for the host exposure model trained on PCR data, based on different optimal thresholding methods.`r`n`r`n(1) pcr_known: The list of observed ('known') host genera based on PCR data.`r`n`r`n(2) pcr_unknown_rs0.8: The list of unobserved ('unknown') predicted host genera when applying an 80% sensitivity threshold.`r`n`r`n(3) pcr_unknown_rs0.85:  The list of unobserved ('unknown') predicted host genera when applying an 85% sensitivity threshold.`r`n`r`n(4) pcr_unknown_rs0.9: The list of unobserved ('unknown') predicted host genera when applying an 90% sensitivity threshold.`r`n`r`n(5) pcr_unknown_rs0.95: The list of unobserved ('unknown') predicted host genera when applying an 95% sensitivity threshold.`r`n`r`n(6) pcr_unknown_mss3: The list of unobserved ('unknown') predicted host genera when applying a threshold that maximizes the sum of sensitivity and specificity, otherwise known as the Youden Index."

$shp.TextFrame.Characters().Text = $readmeText

# ------------------------------------------------------------------
# 2. Rename the data sheets to include their new numbered prefixes.
# ------------------------------------------------------------------
$wb.Worksheets.Item("pcr_known").Name = "(1) pcr_known"
$wb.Worksheets.Item("pcr_unknown_rs0.8").Name = "(2) pcr_unknown_rs0.8"
$wb.Worksheets.Item("pcr_unknown_rs0.85").Name = "(3) pcr_unknown_rs0.85"
$wb.Worksheets.Item("pcr_unknown_rs0.9").Name = "(4) pcr_unknown_rs0.9"
$wb.Worksheets.Item("pcr_unknown_rs0.95").Name = "(5) pcr_unknown_rs0.95"
$wb.Worksheets.Item("pcr_unknown_mss3").Name = "(6) pcr_unknown_mss3"

# ------------------------------------------------------------------
# 3. The "(2) pcr_unknown_rs0.8" sheet used to carry a second,
#    duplicate little table in columns E:G - remove it.
# ------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("(2) pcr_unknown_rs0.8")
$ws2.Range("E1:G70").ClearContents()

# Update the view so it isn't the active tab anymore and scrolls to
# where the old selection used to be.
$ws2.Activate()
$excel.ActiveWindow.ScrollRow = 47
$ws2.Range("E2:G59").Select()

$readme.Activate()
